# Sprint One / Entity Diagram.xlsx
# "Define maybe condition in Entity diagram"
#
# The Condition entity's Condition.ID column previously documented its
# value domain as "-1,0,1". Replace that with the new domain that adds
# the "maybe" state: "no, unknown, maybe, yes".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("K5").Value = "no, unknown, maybe, yes"

# Leave the selection on the edited cell, as in the authored commit.
$ws.Range("K5").Select()
